$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before F - shifts F..K to G..L
$ws.Columns("F").Insert()

# New header + hyperlink formula text in the freshly inserted column
$ws.Range("F3").Value = "Email"
$ws.Range("F4").Value = '&=HYPERLINK("mailto:{{item.Email}}","{{item.Email}}")'

# Update selection to reflect where the user was working
$ws.Range("F4").Select()

# Style the hyperlink formula cell: underline + blue font (matches existing border)
$ws.Range("F4").Font.Underline = 2
$ws.Range("F4").Font.Color = 16711680

# Apply the built-in Hyperlink cell style below it
$ws.Range("F5").Style = "Hyperlink"

# Fix up defined names so they point at the shifted columns
$wb.Names.Item("dates").RefersTo = "=Лист1!`$H`$3"
$wb.Names.Item("PlanData").RefersTo = "=Лист1!`$A`$4:`$I`$5"
$wb.Names.Item("PlanData_Hours").RefersTo = "=Лист1!`$H`$4"
